# "Favicon en logo klaar" - mark the favicon/logo user story as done.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new (20th) userstory row with the favicon/logo story.
$ws.Range("C22").Value = "Als UX Designer wil ik dat er een favicon komt net als de logo"
$ws.Range("D22").Value = "Zodat mensen de website kunnen herkennen."
$ws.Range("E22").Value = "Must H."
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = "Done"

$ws.Range("A22").Value = "Nizam Sarwar"
$ws.Range("A22").HorizontalAlignment = -4108   # xlCenter

# The PDO-connectie userstory (row 21) is already "Done" - mark it green,
# matching the rest of the completed rows in the tracker grid below.
$ws.Range("E27").Interior.Color = 5287936   # RGB(0,176,80) green - "Afgerond"
$ws.Range("F27").Interior.Color = 49407     # RGB(255,192,0) amber - in progress
$ws.Range("F28").Interior.Color = 49407
$ws.Range("F29").Interior.Color = 49407
$ws.Range("F30").Interior.Color = 49407

# Update the selection to reflect where the author left off.
$ws.Range("C22").Select()
